$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front (shifts existing Code/Description/... data
# from A:E to B:F, reusing all existing shared strings untouched).
$ws.Columns.Item(1).Insert()

# Fill in the new "Version" column.
$rng = $ws.Range("A1:A23")

$ws.Range("A1").Formula = "'Version"
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("A$r").Formula = "'2011"
}

# Drop the quote-prefix formatting picked up from the Formula assignments
# above so the cells stay plain (no extra style applied).
$rng.ClearFormats()
